$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Columns B:F for rows 2:25
$arr_BF = New-Object 'object[,]' 24,5
$arr_BF[0,0] = 24.79234167118251
$arr_BF[0,1] = 12.1986898178024
$arr_BF[0,2] = 4.116472192994415
$arr_BF[0,3] = 9.434954983491362
$arr_BF[0,4] = 51.05152519456281
$arr_BF[1,0] = 24.52690093257728
$arr_BF[1,1] = 11.84615480778898
$arr_BF[1,2] = 4.088055010942537
$arr_BF[1,3] = 9.422241046720382
$arr_BF[1,4] = 51.03973662063927
$arr_BF[2,0] = 24.36889176902127
$arr_BF[2,1] = 11.62822034687801
$arr_BF[2,2] = 4.070077268847696
$arr_BF[2,3] = 9.414262293700673
$arr_BF[2,4] = 51.04436539992761
$arr_BF[3,0] = 24.30582228839021
$arr_BF[3,1] = 11.53919234269348
$arr_BF[3,2] = 4.06261818215492
$arr_BF[3,3] = 9.410967310821276
$arr_BF[3,4] = 51.04923214866234
$arr_BF[4,0] = 24.2954314019402
$arr_BF[4,1] = 11.52440077982387
$arr_BF[4,2] = 4.061371576722823
$arr_BF[4,3] = 9.41041754250789
$arr_BF[4,4] = 51.05022014332374
$arr_BF[5,0] = 24.36803575459909
$arr_BF[5,1] = 11.62702035814398
$arr_BF[5,2] = 4.069977211353107
$arr_BF[5,3] = 9.41421803314674
$arr_BF[5,4] = 51.04441897298071
$arr_BF[6,0] = 24.69982855719716
$arr_BF[6,1] = 12.077531711011
$arr_BF[6,2] = 4.106783915898111
$arr_BF[6,3] = 9.430606842699779
$arr_BF[6,4] = 51.04499585664011
$arr_BF[7,0] = 25.38672562499639
$arr_BF[7,1] = 12.94295895982842
$arr_BF[7,2] = 4.174756801362074
$arr_BF[7,3] = 9.461389392392933
$arr_BF[7,4] = 51.14037594266033
$arr_BF[8,0] = 25.90907205686663
$arr_BF[8,1] = 13.56001811907501
$arr_BF[8,2] = 4.222131006246318
$arr_BF[8,3] = 9.483196077570341
$arr_BF[8,4] = 51.2679280290238
$arr_BF[9,0] = 26.14958000972739
$arr_BF[9,1] = 13.83521420947753
$arr_BF[9,2] = 4.24312543220248
$arr_BF[9,3] = 9.49294452965988
$arr_BF[9,4] = 51.33840090837143
$arr_BF[10,0] = 26.24098889567952
$arr_BF[10,1] = 13.93852551327064
$arr_BF[10,2] = 4.250995370109342
$arr_BF[10,3] = 9.496611614127664
$arr_BF[10,4] = 51.36687090373211
$arr_BF[11,0] = 26.22128880977844
$arr_BF[11,1] = 13.91631718487985
$arr_BF[11,2] = 4.249304007952095
$arr_BF[11,3] = 9.495822925471744
$arr_BF[11,4] = 51.36066018860526
$arr_BF[12,0] = 26.15709405171152
$arr_BF[12,1] = 13.84373229952811
$arr_BF[12,2] = 4.243774505067806
$arr_BF[12,3] = 9.493246711512329
$arr_BF[12,4] = 51.34070745429968
$arr_BF[13,0] = 26.11781395306985
$arr_BF[13,1] = 13.79915180377471
$arr_BF[13,2] = 4.240377075355176
$arr_BF[13,3] = 9.491665529449582
$arr_BF[13,4] = 51.32871785420978
$arr_BF[14,0] = 25.89340544380627
$arr_BF[14,1] = 13.54191316125436
$arr_BF[14,2] = 4.220747759447923
$arr_BF[14,3] = 9.482555560874921
$arr_BF[14,4] = 51.2635724360712
$arr_BF[15,0] = 25.75641992205332
$arr_BF[15,1] = 13.38261528453154
$arr_BF[15,2] = 4.208563178644489
$arr_BF[15,3] = 9.476923118327205
$arr_BF[15,4] = 51.2267921370373
$arr_BF[16,0] = 25.67790637174979
$arr_BF[16,1] = 13.29047876752248
$arr_BF[16,2] = 4.201502576353736
$arr_BF[16,3] = 9.473667329097438
$arr_BF[16,4] = 51.20680934567948
$arr_BF[17,0] = 25.65137309081008
$arr_BF[17,1] = 13.25919833246323
$arr_BF[17,2] = 4.199102991963062
$arr_BF[17,3] = 9.472562193019277
$arr_BF[17,4] = 51.20024502459126
$arr_BF[18,0] = 25.77097420122383
$arr_BF[18,1] = 13.39962667813201
$arr_BF[18,2] = 4.209865667849406
$arr_BF[18,3] = 9.477524372024536
$arr_BF[18,4] = 51.23058618038468
$arr_BF[19,0] = 26.17594119980157
$arr_BF[19,1] = 13.86507745832297
$arr_BF[19,2] = 4.245400832118922
$arr_BF[19,3] = 9.49400406961545
$arr_BF[19,4] = 51.34651971668765
$arr_BF[20,0] = 26.44251841378005
$arr_BF[20,1] = 14.16399078580403
$arr_BF[20,2] = 4.268157588183336
$arr_BF[20,3] = 9.504632135079873
$arr_BF[20,4] = 51.43267980058818
$arr_BF[21,0] = 26.30009342809117
$arr_BF[21,1] = 14.00497167792428
$arr_BF[21,2] = 4.256054721202728
$arr_BF[21,3] = 9.498972690635304
$arr_BF[21,4] = 51.38574648995924
$arr_BF[22,0] = 25.76439345387783
$arr_BF[22,1] = 13.39193754458105
$arr_BF[22,2] = 4.209276985439427
$arr_BF[22,3] = 9.477252600152971
$arr_BF[22,4] = 51.22886727116553
$arr_BF[23,0] = 25.19747164687344
$arr_BF[23,1] = 12.7116096789375
$arr_BF[23,2] = 4.156818706975699
$arr_BF[23,3] = 9.453205252937344
$arr_BF[23,4] = 51.10447357438098
$ws.Range("B2:F25").Value = $arr_BF

# Columns I:J for rows 2:25
$arr_IJ = New-Object 'object[,]' 24,2
$arr_IJ[0,0] = 35.71718101469816
$arr_IJ[0,1] = 9.319292756766067
$arr_IJ[1,0] = 35.77477713447178
$arr_IJ[1,1] = 9.32599649678642
$arr_IJ[2,0] = 35.81731954102613
$arr_IJ[2,1] = 9.330338538789693
$arr_IJ[3,0] = 35.83645571585502
$arr_IJ[3,1] = 9.332164926048929
$arr_IJ[4,0] = 35.83974180661267
$arr_IJ[4,1] = 9.332471641837218
$arr_IJ[5,0] = 35.81757033817465
$arr_IJ[5,1] = 9.330362939204598
$arr_IJ[6,0] = 35.73554730859858
$arr_IJ[6,1] = 9.321557413917251
$arr_IJ[7,0] = 35.63187536417959
$arr_IJ[7,1] = 9.306075025469003
$arr_IJ[8,0] = 35.59084323718618
$arr_IJ[8,1] = 9.295778151381864
$arr_IJ[9,0] = 35.57985725244973
$arr_IJ[9,1] = 9.291325737853381
$arr_IJ[10,0] = 35.57680480177332
$arr_IJ[10,1] = 9.289672876170819
$arr_IJ[11,0] = 35.5774128879589
$arr_IJ[11,1] = 9.290027376535662
$arr_IJ[12,0] = 35.57958390798576
$arr_IJ[12,1] = 9.291189091966846
$arr_IJ[13,0] = 35.5810580674183
$arr_IJ[13,1] = 9.291904991719315
$arr_IJ[14,0] = 35.59171603227406
$arr_IJ[14,1] = 9.296073776727901
$arr_IJ[15,0] = 35.60022386693724
$arr_IJ[15,1] = 9.298690427746553
$arr_IJ[16,0] = 35.60584003177402
$arr_IJ[16,1] = 9.300217272223707
$arr_IJ[17,0] = 35.60786557935936
$arr_IJ[17,1] = 9.300737986903322
$arr_IJ[18,0] = 35.59924337275156
$arr_IJ[18,1] = 9.298409623965298
$arr_IJ[19,0] = 35.5789161396868
$arr_IJ[19,1] = 9.290846968894437
$arr_IJ[20,0] = 35.57208886428301
$arr_IJ[20,1] = 9.286097614192455
$arr_IJ[21,0] = 35.57514084371246
$arr_IJ[21,1] = 9.28861479700336
$arr_IJ[22,0] = 35.59968439638746
$arr_IJ[22,1] = 9.298536505269061
$arr_IJ[23,0] = 35.65376970399402
$arr_IJ[23,1] = 9.310073369184019
$ws.Range("I2:J25").Value = $arr_IJ

# Columns L:N for rows 2:25
$arr_LN = New-Object 'object[,]' 24,3
$arr_LN[0,0] = 12.18213002633322
$arr_LN[0,1] = 20.71384453939665
$arr_LN[0,2] = 22.07692054025621
$arr_LN[1,0] = 12.19652173816174
$arr_LN[1,1] = 20.67533557885177
$arr_LN[1,2] = 22.14782543832314
$arr_LN[2,0] = 12.20699692142309
$arr_LN[2,1] = 20.65563771639017
$arr_LN[2,2] = 22.19336377595812
$arr_LN[3,0] = 12.21167793815653
$arr_LN[3,1] = 20.64860851326575
$arr_LN[3,2] = 22.21242588281557
$arr_LN[4,0] = 12.21248012534246
$arr_LN[4,1] = 20.6475017284529
$arr_LN[4,2] = 22.21562166299561
$arr_LN[5,0] = 12.20705838160111
$arr_LN[5,1] = 20.65553887141209
$arr_LN[5,2] = 22.19361880830853
$arr_LN[6,0] = 12.18675235847357
$arr_LN[6,1] = 20.69975094396977
$arr_LN[6,2] = 22.10095374863319
$arr_LN[7,0] = 12.15992201977799
$arr_LN[7,1] = 20.81749914017274
$arr_LN[7,2] = 21.93506412932311
$arr_LN[8,0] = 12.14811119368473
$arr_LN[8,1] = 20.92253340934225
$arr_LN[8,2] = 21.82274419021714
$arr_LN[9,0] = 12.14444908991998
$arr_LN[9,1] = 20.97423735152848
$arr_LN[9,2] = 21.77370371715278
$arr_LN[10,0] = 12.1433078150889
$arr_LN[10,1] = 20.99437066254043
$arr_LN[10,2] = 21.7554274260726
$arr_LN[11,0] = 12.14354269989296
$arr_LN[11,1] = 20.99001011472365
$arr_LN[11,2] = 21.75935048269478
$arr_LN[12,0] = 12.14435027984088
$arr_LN[12,1] = 20.9758826850141
$arr_LN[12,2] = 21.77219422458525
$arr_LN[13,0] = 12.14487689890055
$arr_LN[13,1] = 20.96730108383418
$arr_LN[13,2] = 21.78009967629991
$arr_LN[14,0] = 12.14838490210642
$arr_LN[14,1] = 20.91923254204878
$arr_LN[14,2] = 21.82599034189386
$arr_LN[15,0] = 12.15097475550447
$arr_LN[15,1] = 20.8907420021329
$arr_LN[15,2] = 21.85466814171421
$arr_LN[16,0] = 12.15262542457676
$arr_LN[16,1] = 20.87472492445545
$arr_LN[16,2] = 21.87135629426986
$arr_LN[17,0] = 12.153211988474
$arr_LN[17,1] = 20.86936564623069
$arr_LN[17,2] = 21.87703987497229
$arr_LN[18,0] = 12.1506823962629
$arr_LN[18,1] = 20.89373665312898
$arr_LN[18,2] = 21.85159532913437
$arr_LN[19,0] = 12.14410641605627
$arr_LN[19,1] = 20.98001729739529
$arr_LN[19,2] = 21.76841372889895
$arr_LN[20,0] = 12.14123922095873
$arr_LN[20,1] = 21.03963156551484
$arr_LN[20,2] = 21.71576446347408
$arr_LN[21,0] = 12.14263878373348
$arr_LN[21,1] = 21.00752272035579
$arr_LN[21,2] = 21.7437078465524
$arr_LN[22,0] = 12.15081406807657
$arr_LN[22,1] = 20.89238164145942
$arr_LN[22,2] = 21.85298392178023
$arr_LN[23,0] = 12.165791269619
$arr_LN[23,1] = 20.78235978823906
$arr_LN[23,2] = 21.97825661289527
$ws.Range("L2:N25").Value = $arr_LN
